$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 132.33333
$ws.Range("I6").Value = 113.53333
$ws.Range("J6").Value = 226.33333
$ws.Range("K6").Value = 340.59999
$ws.Range("L6").Value = 678.99999
$ws.Range("M6").Value = -228.59999
$ws.Range("N6").Value = -902.99999
$ws.Range("H8").Value = 40.25
$ws.Range("I8").Value = 40.25
$ws.Range("K8").Value = 120.75
$ws.Range("M8").Value = 18.25
$ws.Range("H33").Value = 31627
$ws.Range("I33").Value = 42021.332
$ws.Range("J33").Value = 444
$ws.Range("K33").Value = 42021.332
$ws.Range("L33").Value = 444
$ws.Range("M33").Value = -41792.332
$ws.Range("N33").Value = -902
$ws.Range("H40").Value = 3769.7
$ws.Range("J40").Value = 2919.4
$ws.Range("L40").Value = 2919.4
$ws.Range("N40").Value = -3269.4
$ws.Range("H112").Value = 1409.2
$ws.Range("I112").Value = 1113
$ws.Range("J112").Value = 1424.7894
$ws.Range("K112").Value = 3339
$ws.Range("L112").Value = 4274.3682
$ws.Range("M112").Value = -2231
$ws.Range("N112").Value = -6490.3682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3478.8235
$ws.Range("I2").Value = 2437.7144
$ws.Range("K2").Value = 2437.7144
$ws.Range("M2").Value = -2324.7144
$ws.Range("H45").Value = 2036.8572
$ws.Range("I45").Value = 1935.75
$ws.Range("K45").Value = 1935.75
$ws.Range("M45").Value = -1558.75
$ws.Range("H110").Value = 9592.538
$ws.Range("I110").Value = 9334.362999999999
$ws.Range("K110").Value = 9334.362999999999
$ws.Range("M110").Value = -7289.362999999999
$ws.Range("H116").Value = 3478.8235
$ws.Range("I116").Value = 2437.7144
$ws.Range("K116").Value = 2437.7144
$ws.Range("M116").Value = -143.7143999999998
$ws.Range("H122").Value = 1689.8572
$ws.Range("I122").Value = 1763.4445
$ws.Range("J122").Value = 1557.4
$ws.Range("K122").Value = 5290.333500000001
$ws.Range("L122").Value = 4672.200000000001
$ws.Range("M122").Value = -2840.333500000001
$ws.Range("N122").Value = -9572.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3478.8235
$ws.Range("I3").Value = 2437.7144
$ws.Range("K3").Value = 2437.7144
$ws.Range("M3").Value = -2323.7144
$ws.Range("H99").Value = 4598.75
$ws.Range("I99").Value = 4467
$ws.Range("K99").Value = 4467
$ws.Range("M99").Value = -2969

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 997.5
$ws.Range("I16").Value = 997.5
$ws.Range("K16").Value = 997.5
$ws.Range("M16").Value = -710.5
$ws.Range("H31").Value = 2449.9412
$ws.Range("J31").Value = 4985.8
$ws.Range("L31").Value = 4985.8
$ws.Range("N31").Value = -5575.8
$ws.Range("H34").Value = 2449.9412
$ws.Range("J34").Value = 4985.8
$ws.Range("L34").Value = 4985.8
$ws.Range("N34").Value = -5389.8
$ws.Range("H50").Value = 44950
$ws.Range("I50").Value = 44950
$ws.Range("K50").Value = 44950
$ws.Range("M50").Value = -44325
$ws.Range("H58").Value = 3318.6
$ws.Range("J58").Value = 3897
$ws.Range("L58").Value = 3897
$ws.Range("N58").Value = -4303
$ws.Range("H60").Value = 38913.043
$ws.Range("H62").Value = 14064.521
$ws.Range("I62").Value = 10055.462
$ws.Range("K62").Value = 10055.462
$ws.Range("M62").Value = -9431.462
$ws.Range("H65").Value = 14064.521
$ws.Range("I65").Value = 10055.462
$ws.Range("K65").Value = 50277.31
$ws.Range("M65").Value = -47157.31
$ws.Range("H99").Value = 3475.7856
$ws.Range("I99").Value = 2862.5
$ws.Range("K99").Value = 2862.5
$ws.Range("M99").Value = -1364.5
$ws.Range("H107").Value = 1774.6957
$ws.Range("I107").Value = 758.4
$ws.Range("K107").Value = 758.4
$ws.Range("M107").Value = 1161.6
$ws.Range("H113").Value = 997.5
$ws.Range("I113").Value = 997.5
$ws.Range("K113").Value = 997.5
$ws.Range("M113").Value = 1172.5
$ws.Range("H126").Value = 3475.7856
$ws.Range("I126").Value = 2862.5
$ws.Range("K126").Value = 8587.5
$ws.Range("M126").Value = -6117.5
$ws.Range("H136").Value = 3318.6
$ws.Range("J136").Value = 3897
$ws.Range("L136").Value = 11691
$ws.Range("N136").Value = -16791

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 100970.1
$ws.Range("I4").Value = 955.05
$ws.Range("J4").Value = 301000.2
$ws.Range("K4").Value = 2865.15
$ws.Range("L4").Value = 903000.6000000001
$ws.Range("M4").Value = -2753.15
$ws.Range("N4").Value = -903224.6000000001
$ws.Range("H26").Value = 962.8276
$ws.Range("I26").Value = 880
$ws.Range("J26").Value = 980.0833
$ws.Range("K26").Value = 2640
$ws.Range("L26").Value = 2940.2499
$ws.Range("M26").Value = -2352
$ws.Range("N26").Value = -3516.2499
$ws.Range("H57").Value = 10525
$ws.Range("I57").Value = 50
$ws.Range("K57").Value = 150
$ws.Range("M57").Value = 409
$ws.Range("H124").Value = 3000
$ws.Range("J124").Value = 3000
$ws.Range("L124").Value = 9000
$ws.Range("N124").Value = -18820
$ws.Range("H131").Value = 506419.88
$ws.Range("I131").Value = 1060.2
$ws.Range("J131").Value = 587929.5
$ws.Range("K131").Value = 3180.6
$ws.Range("L131").Value = 1763788.5
$ws.Range("M131").Value = 1859.4
$ws.Range("N131").Value = -1773868.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1223.25
$ws.Range("I102").Value = 1223.25
$ws.Range("K102").Value = 1223.25
$ws.Range("M102").Value = 398.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5833.1665
$ws.Range("I46").Value = 5999.5
$ws.Range("K46").Value = 5999.5
$ws.Range("M46").Value = -5811.5
$ws.Range("H68").Value = 3011
$ws.Range("I68").Value = 2931.5715
$ws.Range("J68").Value = 3196.3333
$ws.Range("K68").Value = 2931.5715
$ws.Range("L68").Value = 3196.3333
$ws.Range("M68").Value = -2182.5715
$ws.Range("N68").Value = -4694.3333
$ws.Range("H71").Value = 3011
$ws.Range("I71").Value = 2931.5715
$ws.Range("J71").Value = 3196.3333
$ws.Range("K71").Value = 14657.8575
$ws.Range("L71").Value = 15981.6665
$ws.Range("M71").Value = -10913.8575
$ws.Range("N71").Value = -23469.6665
$ws.Range("H93").Value = 1216.48
$ws.Range("I93").Value = 1598.875
$ws.Range("K93").Value = 1598.875
$ws.Range("M93").Value = -350.875
$ws.Range("H122").Value = 3462.1333
$ws.Range("J122").Value = 3492.4614
$ws.Range("L122").Value = 10477.3842
$ws.Range("N122").Value = -15377.3842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2862.875
$ws.Range("I96").Value = 5751.5
$ws.Range("K96").Value = 5751.5
$ws.Range("M96").Value = -4378.5
$ws.Range("H107").Value = 2407.8572
$ws.Range("I107").Value = 963.25
$ws.Range("K107").Value = 2889.75
$ws.Range("M107").Value = -969.75
$ws.Range("H113").Value = 3833
$ws.Range("J113").Value = 6216.857
$ws.Range("L113").Value = 18650.571
$ws.Range("N113").Value = -22990.571
$ws.Range("H126").Value = 3219
$ws.Range("J126").Value = 3937.8
$ws.Range("L126").Value = 11813.4
$ws.Range("N126").Value = -16753.4

Write-Host "Applied all cell updates"